$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.169.28"
$ws.Range("E2").Value = "  -2.16%  "

$ws.Range("D3").Value = "1.902.10"
$ws.Range("E3").Value = "  -2.62%  "

$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'332.40"
$ws.Range("E5").Value = "  -2.98%  "

$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("D7").Value = "'0.4610"
$ws.Range("E7").Value = "  -3.41%  "

$ws.Range("D8").Value = "'0.4124"
$ws.Range("E8").Value = "  -0.82%  "

$ws.Range("D9").Value = "'47.89"
$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("D10").Value = "'0.08023"
$ws.Range("E10").Value = "  -2.94%  "

$ws.Range("D11").Value = "'1.013"
$ws.Range("E11").Value = "  -2.39%  "

$ws.Range("D12").Value = "'22.17"
$ws.Range("E12").Value = "  -2.58%  "

$ws.Range("D13").Value = "1.894.23"
$ws.Range("E13").Value = "  -3.10%  "

$ws.Range("D14").Value = "'5.949"
$ws.Range("E14").Value = "  -3.85%  "

$ws.Range("D15").Value = "'7.120"
$ws.Range("E15").Value = "  -3.93%  "

$ws.Range("D16").Value = "'89.26"
$ws.Range("E16").Value = "  -3.23%  "

$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("D18").Value = "'0.00001030"
$ws.Range("E18").Value = "  -2.50%  "

$ws.Range("D19").Value = "'0.06578"
$ws.Range("E19").Value = "  -1.79%  "

$ws.Range("D20").Value = "'17.65"
$ws.Range("E20").Value = "  -2.22%  "

$ws.Range("D22").Value = "29.147.63"
$ws.Range("E22").Value = "  -2.11%  "

$ws.Range("D23").Value = "'5.504"
$ws.Range("E23").Value = "  -1.51%  "

$ws.Range("D24").Value = "'11.42"
$ws.Range("E24").Value = "  +1.15%  "

$ws.Range("D25").Value = "'2.194"
$ws.Range("E25").Value = "  -3.36%  "

$ws.Range("D26").Value = "2.125.68"
$ws.Range("E26").Value = "  -2.61%  "

$ws.Range("D27").Value = "'157.25"
$ws.Range("E27").Value = "  -2.98%  "

$ws.Range("D28").Value = "'19.76"
$ws.Range("E28").Value = "  -2.19%  "

$ws.Range("D29").Value = "'2.126"
$ws.Range("E29").Value = "  -2.68%  "

$ws.Range("D30").Value = "'5.644"
$ws.Range("E30").Value = "  -1.15%  "

$ws.Range("D31").Value = "'117.24"
$ws.Range("E31").Value = "  -4.70%  "

$ws.Range("E32").Value = "  +3.45%  "

$ws.Range("D33").Value = "'0.09427"
$ws.Range("E33").Value = "  -2.17%  "

$ws.Range("D34").Value = "'1.427"
$ws.Range("E34").Value = "  -3.65%  "

$ws.Range("D35").Value = "'3.544"
$ws.Range("E35").Value = "  -3.92%  "

$ws.Range("D36").Value = "'5.359"
$ws.Range("E36").Value = "  -3.03%  "

$ws.Range("D37").Value = "'0.06103"
$ws.Range("E37").Value = "  -3.12%  "

$ws.Range("D38").Value = "'0.02247"
$ws.Range("E38").Value = "  -3.25%  "

$ws.Range("D39").Value = "'8.446"
$ws.Range("E39").Value = "  -0.62%  "

$ws.Range("D40").Value = "'1.177"
$ws.Range("E40").Value = "  -0.92%  "

$ws.Range("D41").Value = "'0.5848"
$ws.Range("E41").Value = "  -4.26%  "

$ws.Range("D42").Value = "'0.9999"
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("E43").Value = "  -3.66%  "

$ws.Range("D44").Value = "'10.16"
$ws.Range("E44").Value = "  -5.41%  "

$ws.Range("D47").Value = "'0.07495"
$ws.Range("E47").Value = "  +1.78%  "

$ws.Range("D48").Value = "'0.5554"
$ws.Range("E48").Value = "  -3.03%  "

$ws.Range("D49").Value = "'12.09"
$ws.Range("E49").Value = "  -3.74%  "

$ws.Range("D50").Value = "'1.924"
$ws.Range("E50").Value = "  -3.33%  "

$ws.Range("D51").Value = "'113.00"
$ws.Range("E51").Value = "  -0.47%  "

# Row 45 and 46 swap (RenderToken <-> WEMIXTOKEN)
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.256"
$ws.Range("E45").Value = "  -1.35%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'2.346"
$ws.Range("E46").Value = "  -2.46%  "
